# Correção das notas do fórum para matc65 em 2021.2
# Zera as colunas de visualizações diárias (B:J) para as linhas de alunos
# (linhas 2 a 50), já que a linha 1 é o cabeçalho.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:J50").Value = 0
